$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates per diff (B/C = coin name/link swaps, D = price, E = volume % change).
# Numeric-looking D-column values are written with a leading apostrophe so Excel stores them
# as literal text (matching the original inlineStr cells) instead of re-parsing them as numbers.

$ws.Range("D2").Value = "70.907.01"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "3.847.05"
$ws.Range("E3").Value = "  +1.40%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'707.26"
$ws.Range("E5").Value = "  +1.54%  "
$ws.Range("D6").Value = "'172.73"
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("D7").Value = "3.845.53"
$ws.Range("E7").Value = "  +1.45%  "
$ws.Range("E9").Value = "  -0.66%  "
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("D11").Value = "'7.35"
$ws.Range("E11").Value = "  -1.04%  "
$ws.Range("D12").Value = "'0.458"
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("D13").Value = "'0.0000256"
$ws.Range("E13").Value = "  -0.93%  "
$ws.Range("E14").Value = "  +1.06%  "
$ws.Range("D15").Value = "4.495.83"
$ws.Range("E15").Value = "  +1.44%  "
$ws.Range("D16").Value = "3.844.26"
$ws.Range("E16").Value = "  +1.29%  "
$ws.Range("D17").Value = "70.967.35"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("E19").Value = "  +0.74%  "
$ws.Range("D20").Value = "'17.38"
$ws.Range("E20").Value = "  -2.43%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'492.57"
$ws.Range("E21").Value = "  +1.66%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'10.65"
$ws.Range("E22").Value = "  -3.87%  "
$ws.Range("E23").Value = "  +0.47%  "
$ws.Range("D24").Value = "'85.10"
$ws.Range("E24").Value = "  +1.15%  "
$ws.Range("E25").Value = "  +2.13%  "
$ws.Range("E26").Value = "  +1.27%  "
$ws.Range("D27").Value = "'12.14"
$ws.Range("E27").Value = "  -1.95%  "
$ws.Range("D28").Value = "'2.11"
$ws.Range("E28").Value = "  -2.01%  "
$ws.Range("D29").Value = "'3.19"
$ws.Range("E29").Value = "  +4.17%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("D31").Value = "'7.49"
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("E32").Value = "  -0.58%  "
$ws.Range("D33").Value = "'29.45"
$ws.Range("E33").Value = "  -0.44%  "
$ws.Range("E34").Value = "  +0.25%  "
$ws.Range("D35").Value = "'9.17"
$ws.Range("E35").Value = "  -0.37%  "
$ws.Range("D36").Value = "3.803.15"
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("E38").Value = "  +0.77%  "
$ws.Range("E39").Value = "  +6.80%  "
$ws.Range("D40").Value = "'6.06"
$ws.Range("E40").Value = "  +1.41%  "
$ws.Range("D41").Value = "'1.04"
$ws.Range("E41").Value = "  +6.53%  "
$ws.Range("E42").Value = "  -4.88%  "
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").Value = "'163.57"
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").Value = "'0.000309"
$ws.Range("E46").Value = "  -4.84%  "
$ws.Range("D47").Value = "'48.74"
$ws.Range("E47").Value = "  -1.03%  "
$ws.Range("E48").Value = "  +0.78%  "
$ws.Range("D49").Value = "'413.73"
$ws.Range("E49").Value = "  +2.76%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "'8.62"
$ws.Range("E50").Value = "  +0.97%  "
$ws.Range("B51").Value = "TheGraph"
$ws.Range("C51").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D51").Value = "'0.298"
$ws.Range("E51").Value = "  -0.47%  "
